$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4 and B5 to use the correct "secret_sauce" value (fixing typos)
$ws.Range("B4").Value = "secret_sauce"
$ws.Range("B5").Value = "secret_sauce"

# Update the selection to match the new state (B5 selected)
$ws.Range("B5").Select()
